$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "AurXX1"
$ws.Range("B2").Value = "Bert"
$ws.Range("C2").Value = 36891
$ws.Range("D2").Value = 0.625
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 36892
$ws.Range("G2").Value = 0.708333333333333

# Row 3
$ws.Range("A3").Value = "AurXX1"
$ws.Range("B3").Value = "Bert"
$ws.Range("C3").Value = 36892
$ws.Range("D3").Value = 0.625
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 36892
$ws.Range("G3").Value = 0.708333333333333

# Row 4 (new)
$ws.Range("A4").Value = "AurXX2"
$ws.Range("B4").Value = "Bert"
$ws.Range("C4").Value = 36893
$ws.Range("D4").Value = 0.666666666666667
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 36893
$ws.Range("G4").Value = 0.708333333333333

# Row 5 (previously row 4)
$ws.Range("A5").Value = "AurXX1"
$ws.Range("B5").Value = "Ernie"
$ws.Range("C5").Value = 36896
$ws.Range("D5").Value = 0.0416666666666667
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = 36896
$ws.Range("G5").Value = 0.541666666666667

# Row 6 (new)
$ws.Range("A6").Value = "AurXX1"
$ws.Range("B6").Value = "Ernie"
$ws.Range("C6").Value = 36923
$ws.Range("D6").Value = 0.0416666666666667
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 36923
$ws.Range("G6").Value = 0.541666666666667

$null = $ws.Range("A7:G9").Select()

# Number format for the dateStart/dateEnd columns changes from DD/MM/YYYY to DD/MM/YY
$ws.Range("C1:C6").NumberFormat = "DD/MM/YY"
$ws.Range("F1:F6").NumberFormat = "DD/MM/YY"

# Column widths shrink slightly (closest value achievable given pixel snapping)
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 10.5
$ws.Columns.Item(3).ColumnWidth = 10.5
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 10.5
$ws.Columns.Item(6).ColumnWidth = 10.5
$ws.Columns.Item(7).ColumnWidth = 10.5

# Tab ratio (book view) shifts slightly
$excel.ActiveWindow.TabRatio = 985
